# Apply cell value updates for the cryptos worksheet per the Dec 2 2023 data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "38.781.15"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.66%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.092.65"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.15%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "229.09"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.16%  "
$ws.Range("E6").Value = "  +0.56%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "61.22"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.36%  "
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("E9").Value = "  +1.86%  "
$ws.Range("E10").Value = "  +0.08%  "
$ws.Range("E11").Value = "  -0.20%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "15.35"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +4.65%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.403.49"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.05%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "22.11"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.16%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.804"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +3.72%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.50"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.14%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.083.88"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.37%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "38.673.82"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.49%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "71.88"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.42%  "
$ws.Range("E20").Value = "  +0.76%  "
$ws.Range("E21").Value = "  +0.53%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "226.25"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.00%  "
$ws.Range("E24").Value = "  -2.60%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.34"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.83%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "171.42"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.66%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.52"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.10%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.138"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +5.21%  "
$ws.Range("E29").Value = "  +6.38%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "19.54"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.86%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.47"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +3.49%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.120"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.52%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.74"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.23%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0612"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.99%  "
$ws.Range("B36").Value = "LidoDAOToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.40"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.41%  "
$ws.Range("B37").Value = "THORChain"
$ws.Range("C37").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.43"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.58%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.58"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.42%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.999"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.00%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "17.97"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.62%  "
$ws.Range("E41").Value = "  +4.81%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "100.84"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.98%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.533.29"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.75%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.80"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.88%  "
$ws.Range("E45").Value = "  +0.89%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "7.72"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +6.51%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.13"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.42%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "4.12"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.64%  "
$ws.Range("E49").Value = "  +1.64%  "
$ws.Range("E50").Value = "  -1.10%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.289.70"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.11%  "
